# Apply "想去人数" (F column) count updates as captured in commit
# "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 166
$ws.Cells.Item(6, 6).Value = 330
$ws.Cells.Item(7, 6).Value = 5904
$ws.Cells.Item(8, 6).Value = 10001
$ws.Cells.Item(9, 6).Value = 434
$ws.Cells.Item(12, 6).Value = 3971
$ws.Cells.Item(16, 6).Value = 118
$ws.Cells.Item(18, 6).Value = 666
$ws.Cells.Item(19, 6).Value = 3969
$ws.Cells.Item(20, 6).Value = 146
$ws.Cells.Item(21, 6).Value = 92
$ws.Cells.Item(22, 6).Value = 5517
$ws.Cells.Item(23, 6).Value = 445
$ws.Cells.Item(24, 6).Value = 2176
$ws.Cells.Item(25, 6).Value = 142
$ws.Cells.Item(26, 6).Value = 382
$ws.Cells.Item(27, 6).Value = 8242
$ws.Cells.Item(30, 6).Value = 2226
$ws.Cells.Item(31, 6).Value = 2258
$ws.Cells.Item(33, 6).Value = 185
$ws.Cells.Item(34, 6).Value = 1347
$ws.Cells.Item(36, 6).Value = 288
$ws.Cells.Item(38, 6).Value = 264
$ws.Cells.Item(41, 6).Value = 1196
$ws.Cells.Item(42, 6).Value = 1191
$ws.Cells.Item(43, 6).Value = 45
$ws.Cells.Item(44, 6).Value = 66
$ws.Cells.Item(46, 6).Value = 1374

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 617
$ws.Cells.Item(4, 6).Value = 76

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 166
$ws.Cells.Item(4, 6).Value = 617
$ws.Cells.Item(6, 6).Value = 76
$ws.Cells.Item(7, 6).Value = 330
$ws.Cells.Item(8, 6).Value = 5904
$ws.Cells.Item(9, 6).Value = 10001
$ws.Cells.Item(10, 6).Value = 434
$ws.Cells.Item(11, 6).Value = 3971
$ws.Cells.Item(14, 6).Value = 118
$ws.Cells.Item(18, 6).Value = 666
$ws.Cells.Item(19, 6).Value = 3969
$ws.Cells.Item(21, 6).Value = 146
$ws.Cells.Item(22, 6).Value = 5517
$ws.Cells.Item(23, 6).Value = 445
$ws.Cells.Item(24, 6).Value = 2176
$ws.Cells.Item(25, 6).Value = 142
$ws.Cells.Item(26, 6).Value = 382
$ws.Cells.Item(27, 6).Value = 8242
$ws.Cells.Item(30, 6).Value = 2226
$ws.Cells.Item(31, 6).Value = 2258
$ws.Cells.Item(33, 6).Value = 185
$ws.Cells.Item(34, 6).Value = 1347
$ws.Cells.Item(36, 6).Value = 288
$ws.Cells.Item(37, 6).Value = 264
$ws.Cells.Item(40, 6).Value = 1196
$ws.Cells.Item(41, 6).Value = 1191
$ws.Cells.Item(42, 6).Value = 45
$ws.Cells.Item(43, 6).Value = 66
$ws.Cells.Item(45, 6).Value = 1374
